# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Column width for the new column G (matches width=17 in the OOXML) ---
$ws.Columns.Item(7).ColumnWidth = 16.14

# --- Copy formatting from column F into column G so styles (borders, number
#     format, bold header, totals format) line up with the rest of the sheet ---
$ws.Range("F1:F264").Copy()
$ws.Range("G1:G264").PasteSpecial(-4122)

# --- Header ---
$ws.Range("G1").Value = "PRESUPUESTO"

# --- Budget values for rows 2-263 (detail rows), default 0 ---
$budget = @{
    139 = 1000
    141 = 3000
    142 = 5000
    144 = 6000
    146 = 7000
    148 = 6000
    149 = 6000
    150 = 1000
    151 = 400
    152 = 6500
    154 = 4000
    156 = 500
    157 = 4000
}

$data = New-Object 'object[,]' 262,1
for ($r = 2; $r -le 263; $r++) {
    $idx = $r - 2
    if ($budget.ContainsKey($r)) {
        $data[$idx,0] = $budget[$r]
    } else {
        $data[$idx,0] = 0
    }
}
$ws.Range("G2:G263").Value = $data

# --- Totals row ---
$ws.Range("G264").Value = 50400
